$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = 150
    $ws.Cells.Item($r, 2).Value = 26
    $ws.Cells.Item($r, 3).Value = 124
}
